# Sync attendance_reports: fix ordering of the "Recorded By" (column G)
# comma-separated list of recorders. The last recorder in the list should
# be moved to the front (right-rotation of the comma separated values).
#
# e.g. "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#      "System, system, backup@backdoor.com" -> "backup@backdoor.com, System, system"
#      "admin@admin.com, System"             -> "System, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts at row 2; row 1 is the header).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $n = $parts.Count

        if ($n -gt 1) {
            $lastPart = $parts[$n - 1]
            $rest = $parts[0..($n - 2)]
            $newParts = @($lastPart) + $rest
            $cell.Value2 = ($newParts -join ", ")
        }
    }
}
